$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.AddChart2(-1, 4)
$chart = $shp.Chart
$chart.SeriesCollection().NewSeries()
$s = $chart.SeriesCollection(1)
$s.Values = '=Feuil1!$M$2:$M$8'
try {
  $chart.Location(1, "Feuil1")
  Write-Host "Location ok"
} catch { Write-Host "Location err: $_" }
